# "removed fail safe criteria"
#
# The "Failure Actions" block originally listed Classification / Driver
# Notifications / Fail Safe Actions / Power Level / High Voltage State,
# followed by the "Operating Conditions" heading and its "Enable
# Criteria .. Ignition ON" line.
#
# After the edit:
#   - the "Operating Conditions" heading moves up, right after "Failure
#     Actions" (and now carries the _GoBack bookmark that used to sit at
#     the very end of the Repair Actions paragraph);
#   - the "Classification" line is rewritten into a merged
#     "Enable Criteria ... Ignition ON lassification ... Slight Fault"
#     line (a leftover of the manual edit, complete with the spell-check
#     squiggle markup on "lassification");
#   - the "Driver Notifications", "Fail Safe Actions", "Power Level" and
#     "High Voltage State" lines are removed outright;
#   - the old "Enable Criteria / Ignition ON" line is reduced to a single
#     word, "Driver".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The _GoBack bookmark is relocated from the tail of the repair
#    actions paragraph up onto the "Operating Conditions" heading, so
#    drop it from its old spot first (its new spot is (re)created below).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Locate the six paragraphs that make up the old block, from
#    "Classification" through the "Operating Conditions" heading, and
#    replace the whole block in one shot with the two paragraphs that
#    remain afterwards.
# ---------------------------------------------------------------------
$blockStart = $null
$blockEnd = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13)

    if ($blockStart -eq $null -and $text -eq "Classification`tSlight Fault") {
        $blockStart = $para
    }
    if ($para.Style.NameLocal -eq "Heading 1" -and $text -eq "Operating Conditions") {
        $blockEnd = $para
    }
}

$blockRange = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)

$newBlockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="181" w:after="240"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="0071BB"/></w:rPr><w:t>Operating</w:t></w:r><w:r><w:rPr><w:color w:val="0071BB"/><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="0071BB"/></w:rPr><w:t>Conditions</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/><w:tabs><w:tab w:val="left" w:pos="2668"/></w:tabs><w:spacing w:before="52"/></w:pPr><w:r><w:t>Enable</w:t></w:r><w:r><w:rPr><w:spacing w:val="-2"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Criteria</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>Ignition ON</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lassification</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/><w:t>Slight Fault</w:t></w:r></w:p>'

$blockRange.InsertXML($newBlockXml)

# ---------------------------------------------------------------------
# 3. The paragraph that used to read "Enable Criteria .. Ignition ON"
#    (its own spacing/indent settings are untouched) now just says
#    "Driver".
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq "Enable Criteria`tIgnition ON") {
        $target = $para
        break
    }
}

$enableCriteriaRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$enableCriteriaRange.Text = "Driver"
